# Update visitor/attendance numbers ("output generated" figures) across sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1242
$ws1.Range("F5").Value = 4586
$ws1.Range("F6").Value = 1760
$ws1.Range("F7").Value = 6314
$ws1.Range("F9").Value = 1896
$ws1.Range("F10").Value = 505
$ws1.Range("F16").Value = 7812
$ws1.Range("F19").Value = 179
$ws1.Range("F21").Value = 1738
$ws1.Range("F29").Value = 793
$ws1.Range("F32").Value = 8
$ws1.Range("F33").Value = 73

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 89
$ws2.Range("F20").Value = 32

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2269

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2269
$ws4.Range("F5").Value = 1242
$ws4.Range("F10").Value = 4590
$ws4.Range("F12").Value = 1760
$ws4.Range("F13").Value = 6314
$ws4.Range("F15").Value = 1896
$ws4.Range("F17").Value = 505
$ws4.Range("F23").Value = 7812
$ws4.Range("F26").Value = 179
$ws4.Range("F28").Value = 1738
$ws4.Range("F35").Value = 793
$ws4.Range("F37").Value = 359
$ws4.Range("F40").Value = 89
$ws4.Range("F45").Value = 32
